# Applies the "adding medicine code worked out" edit to the Manufacturer
# Packaging Slip / Quick Returns invoice sheet.
#
# Notes on technique:
#   - Excel auto-detects date-like ("mm/dd/yy(yy)") and purely-numeric text
#     typed into a cell via .Value and silently converts it to a real
#     date-serial / number, even though the source workbook stores these
#     NDC codes / lot-style dates as literal text. To avoid that coercion
#     (and to avoid it creating brand-new cell styles via quotePrefix),
#     every *text* write goes through a scratch cell (E33 - blank, already
#     Text-formatted, inside the existing used range so the sheet
#     dimension does not change) and is transplanted with
#     Copy/PasteSpecial(values), which carries the literal string across
#     without re-parsing it.
#   - L34 is a numeric cell whose NumberFormat is Text ("@"). Typing a
#     number into a Text-formatted cell is stored as literal digits
#     (text), not a number, so it is round-tripped through General format
#     momentarily and back to "@" (a no-op for the style since it already
#     is "@") so the stored value stays a genuine number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell: blank, Text-formatted (NumberFormat "@"), already inside
# the sheet's used range (A1:O57) so round-tripping through it leaves no
# trace (no dimension change, no leftover value/format).
$scratch = $ws.Range("E33")

function Set-TextValue($addr, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $dest = $ws.Range($addr)
    $dest.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.ClearContents()
}

function Set-NumberInTextCell($addr, $number) {
    $dest = $ws.Range($addr)
    $fmt = $dest.NumberFormat
    $dest.NumberFormat = "General"
    $dest.Value = $number
    $dest.NumberFormat = $fmt
}

$excel.CutCopyMode = 0

# ---- Header block ----
Set-TextValue "B5"  "02/20/2024"
Set-TextValue "L5"  "QR02202024AP101"

# ---- "Issue Credit to" block ----
Set-TextValue "D10" "KINARY"
Set-TextValue "D11" "152.35 TENTH AVE"
Set-TextValue "D12" "WHITESTONE, NY,11357"
Set-TextValue "D13" "Account#: "
Set-TextValue "D14" "Phone: 718-767-1234/ 888-527-6806"
Set-TextValue "D15" "DEA: RK0416900"

# ---- Line item rows 20-21 (existing Ascend/Aripiprazole rows -> Avet/Rasagiline) ----
Set-TextValue "B20" "2315574603"
Set-TextValue "C20" "Avet Pharmaceuticals Inc."
Set-TextValue "D20" "Rasagiline mesylate"
Set-TextValue "E20" ".5 mg/1"
Set-TextValue "G20" "RCY01AD6"
Set-TextValue "H20" "12/23/31"
Set-TextValue "I20" "30 CT"

Set-TextValue "B21" "2315574603"
Set-TextValue "C21" "Avet Pharmaceuticals Inc."
Set-TextValue "D21" "Rasagiline mesylate"
Set-TextValue "E21" ".5 mg/1"
Set-TextValue "G21" "RCY01AD6"
Set-TextValue "H21" "12/23/31"
Set-TextValue "I21" "30 CT"

# ---- Line item row 22 (previously blank) ----
Set-TextValue "B22" "2315574603"
Set-TextValue "C22" "Avet Pharmaceuticals Inc."
Set-TextValue "D22" "Rasagiline mesylate"
Set-TextValue "E22" ".5 mg/1"
Set-TextValue "F22" "TABLET"
Set-TextValue "G22" "RCY01AD6"
Set-TextValue "H22" "12/23/31"
Set-TextValue "I22" "30 CT"
$ws.Range("J22").Value = 1
$ws.Range("L22").Value = 1

# ---- Line item row 23 (previously blank) ----
Set-TextValue "B23" "6330490190"
Set-TextValue "C23" "Sun Pharmaceutical Industries, Inc."
Set-TextValue "D23" "Fenofibrate"
Set-TextValue "E23" "160 mg/1"
Set-TextValue "F23" "TABLET, FILM COATED"
Set-TextValue "G23" "MHC1672A"
Set-TextValue "H23" "11/23/30"
Set-TextValue "I23" "90 CT"
$ws.Range("J23").Value = 24
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 24

# ---- Total price ----
Set-NumberInTextCell "L34" 27

$scratch.ClearContents()
$excel.CutCopyMode = 0
